# resultado_clasificacion_unificado.xlsx
# "Modelos escalados, implementacion 3er clasificador tipo 0 vs tipo 3, pruebas correctas"
#
# - Adds two new result columns: sospecha_tipo3_% (E) and porcentaje_tipo_plagio (F)
# - Re-runs the comparison of codigo1.py vs codigo4.py through the (now 3) classifiers,
#   appending one result row per classifier run (rows 2-9), replacing the old single
#   codigo4.py/codigo1.py row with codigo1.py/codigo4.py ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------------
# E1/F1 are new headers; give them the same bold/centered/bordered look as the
# other header cells by copying A1:B1's formatting onto them before writing text.
$ws.Range("A1:B1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "sospecha_tipo3_%"
$ws.Range("F1").Value = "porcentaje_tipo_plagio"

# ---- Data rows --------------------------------------------------------------
# file1, file2, plagio_predicho_binario, tipo_plagio_predicho, sospecha_tipo3_%, porcentaje_tipo_plagio
$rows = @(
    @("codigo1.py", "codigo4.py", 1, 0, 0,    $null),
    @("codigo1.py", "codigo4.py", 1, 2, 0,    $null),
    @("codigo1.py", "codigo4.py", 0, 0, 34.92, $null),
    @("codigo1.py", "codigo4.py", 0, 0, -1,    0),
    @("codigo1.py", "codigo4.py", 0, 0, 34.92, 0),
    @("codigo1.py", "codigo4.py", 1, 1, 0,    99),
    @("codigo1.py", "codigo4.py", 1, 2, 0,    85.33),
    @("codigo1.py", "codigo4.py", 0, 0, 31.4,  0)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    if ($row[5] -ne $null) {
        $ws.Cells.Item($r, 6).Value = $row[5]
    } else {
        # no plagiarism-type percentage produced for this run yet -> leave blank
        $ws.Cells.Item($r, 6).Value = ""
    }
    $r++
}
